$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "23+62="
$t.Cell(1,2).Range.Text = "54-50="
$t.Cell(1,3).Range.Text = "57+25="
$t.Cell(1,4).Range.Text = "64+7="
$t.Cell(1,5).Range.Text = "12-3="

$t.Cell(2,1).Range.Text = "88+8="
$t.Cell(2,2).Range.Text = "15-3="
$t.Cell(2,3).Range.Text = "41+4="
$t.Cell(2,4).Range.Text = "67-61="
$t.Cell(2,5).Range.Text = "51-25="

$t.Cell(3,1).Range.Text = "51-0="
$t.Cell(3,2).Range.Text = "38-12="
$t.Cell(3,3).Range.Text = "39+15="
$t.Cell(3,4).Range.Text = "9+34="
$t.Cell(3,5).Range.Text = "69+28="

$t.Cell(4,1).Range.Text = "37+19="
$t.Cell(4,2).Range.Text = "8+65="
$t.Cell(4,3).Range.Text = "16-15="
$t.Cell(4,4).Range.Text = "0+21="
$t.Cell(4,5).Range.Text = "85-45="

$t.Cell(5,1).Range.Text = "5+76="
$t.Cell(5,2).Range.Text = "87-24="
$t.Cell(5,3).Range.Text = "29-22="
$t.Cell(5,4).Range.Text = "4+85="
$t.Cell(5,5).Range.Text = "7+82="

$t.Cell(6,1).Range.Text = "15+0="
$t.Cell(6,2).Range.Text = "66-35="
$t.Cell(6,3).Range.Text = "60-27="
$t.Cell(6,4).Range.Text = "49+35="
$t.Cell(6,5).Range.Text = "93-8="

$t.Cell(7,1).Range.Text = "27-6="
$t.Cell(7,2).Range.Text = "80-60="
$t.Cell(7,3).Range.Text = "40+47="
$t.Cell(7,4).Range.Text = "47+40="
$t.Cell(7,5).Range.Text = "88+11="

$t.Cell(8,1).Range.Text = "25+57="
$t.Cell(8,2).Range.Text = "13+57="
$t.Cell(8,3).Range.Text = "17-2="
$t.Cell(8,4).Range.Text = "85+9="
$t.Cell(8,5).Range.Text = "79-44="

$t.Cell(9,1).Range.Text = "38+28="
$t.Cell(9,2).Range.Text = "42-23="
$t.Cell(9,3).Range.Text = "66-47="
$t.Cell(9,4).Range.Text = "12+48="
$t.Cell(9,5).Range.Text = "70-3="

$t.Cell(10,1).Range.Text = "25+7="
$t.Cell(10,2).Range.Text = "77-54="
$t.Cell(10,3).Range.Text = "93-20="
$t.Cell(10,4).Range.Text = "63-30="
$t.Cell(10,5).Range.Text = "6+28="

$t.Cell(11,1).Range.Text = "65-18="
$t.Cell(11,2).Range.Text = "21+47="
$t.Cell(11,3).Range.Text = "32-21="
$t.Cell(11,4).Range.Text = "91-41="
$t.Cell(11,5).Range.Text = "47+7="

$t.Cell(12,1).Range.Text = "91-36="
$t.Cell(12,2).Range.Text = "20+33="
$t.Cell(12,3).Range.Text = "63-27="
$t.Cell(12,4).Range.Text = "4+4="
$t.Cell(12,5).Range.Text = "93-28="

$t.Cell(13,1).Range.Text = "31+40="
$t.Cell(13,2).Range.Text = "28-3="
$t.Cell(13,3).Range.Text = "60-59="
$t.Cell(13,4).Range.Text = "79-70="
$t.Cell(13,5).Range.Text = "54+41="

$t.Cell(14,1).Range.Text = "86+0="
$t.Cell(14,2).Range.Text = "15-0="
$t.Cell(14,3).Range.Text = "65+31="
$t.Cell(14,4).Range.Text = "40+19="
$t.Cell(14,5).Range.Text = "12-9="

$t.Cell(15,1).Range.Text = "3+28="
$t.Cell(15,2).Range.Text = "22-20="
$t.Cell(15,3).Range.Text = "40+8="
$t.Cell(15,4).Range.Text = "34+56="
$t.Cell(15,5).Range.Text = "41-4="

$t.Cell(16,1).Range.Text = "48-37="
$t.Cell(16,2).Range.Text = "28+49="
$t.Cell(16,3).Range.Text = "14+75="
$t.Cell(16,4).Range.Text = "7+55="
$t.Cell(16,5).Range.Text = "75-51="

$t.Cell(17,1).Range.Text = "53+39="
$t.Cell(17,2).Range.Text = "59-30="
$t.Cell(17,3).Range.Text = "98-2="
$t.Cell(17,4).Range.Text = "48+16="
$t.Cell(17,5).Range.Text = "41+32="

$t.Cell(18,1).Range.Text = "41-8="
$t.Cell(18,2).Range.Text = "74-8="
$t.Cell(18,3).Range.Text = "34+52="
$t.Cell(18,4).Range.Text = "5+2="
$t.Cell(18,5).Range.Text = "96-24="

$t.Cell(19,1).Range.Text = "7+87="
$t.Cell(19,2).Range.Text = "23+28="
$t.Cell(19,3).Range.Text = "57+22="
$t.Cell(19,4).Range.Text = "47-29="
$t.Cell(19,5).Range.Text = "82-75="

$t.Cell(20,1).Range.Text = "19+26="
$t.Cell(20,2).Range.Text = "21+68="
$t.Cell(20,3).Range.Text = "36+37="
$t.Cell(20,4).Range.Text = "11+30="
$t.Cell(20,5).Range.Text = "39-25="
